$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they stay as text (matching source formatting)
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values cell by cell, in diff order
$ws.Range('D2').Value = '29.207.73'
$ws.Range('E2').Value = '  -2.06%  '
$ws.Range('D3').Value = '1.842.31'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('D4').Value = '0.9987'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '240.44'
$ws.Range('E5').Value = '  -2.66%  '
$ws.Range('D6').Value = '0.6852'
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('D7').Value = '0.9993'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '0.3008'
$ws.Range('E8').Value = '  -2.75%  '
$ws.Range('D9').Value = '0.07483'
$ws.Range('E9').Value = '  -3.86%  '
$ws.Range('D10').Value = '23.32'
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('D11').Value = '0.07646'
$ws.Range('E11').Value = '  -2.59%  '
$ws.Range('D12').Value = '1.843.50'
$ws.Range('E12').Value = '  -1.74%  '
$ws.Range('D13').Value = '5.064'
$ws.Range('E13').Value = '  -2.59%  '
$ws.Range('D14').Value = '0.6835'
$ws.Range('E14').Value = '  -1.94%  '
$ws.Range('D15').Value = '88.09'
$ws.Range('E15').Value = '  -5.14%  '
$ws.Range('D16').Value = '6.152'
$ws.Range('E16').Value = '  -7.56%  '
$ws.Range('D17').Value = '29.197.28'
$ws.Range('E17').Value = '  -2.09%  '
$ws.Range('D18').Value = '0.000008212'
$ws.Range('E18').Value = '  -2.49%  '
$ws.Range('D19').Value = '2.080.37'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').Value = '231.11'
$ws.Range('E20').Value = '  -5.32%  '
$ws.Range('D21').Value = '12.56'
$ws.Range('E21').Value = '  -2.23%  '
$ws.Range('D22').Value = '0.9986'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').Value = '7.382'
$ws.Range('E23').Value = '  -3.69%  '
$ws.Range('D24').Value = '0.9991'
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').Value = '160.33'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').Value = '0.1452'
$ws.Range('E26').Value = '  -4.31%  '
$ws.Range('D27').Value = '8.728'
$ws.Range('E27').Value = '  -2.82%  '
$ws.Range('D28').Value = '18.14'
$ws.Range('E28').Value = '  -1.73%  '
$ws.Range('D29').Value = '1.508'
$ws.Range('E29').Value = '  -2.44%  '
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('D31').Value = '4.155'
$ws.Range('E31').Value = '  -2.20%  '
$ws.Range('D32').Value = '1.197'
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('D33').Value = '0.05284'
$ws.Range('E33').Value = '  +2.90%  '
$ws.Range('D34').Value = '0.7574'
$ws.Range('E34').Value = '  -4.12%  '
$ws.Range('D35').Value = '1.867'
$ws.Range('E35').Value = '  -3.64%  '
$ws.Range('D36').Value = '1.136'
$ws.Range('E36').Value = '  -2.69%  '
$ws.Range('D37').Value = '2.684'
$ws.Range('E37').Value = '  -0.92%  '
$ws.Range('D38').Value = '1.310.66'
$ws.Range('E38').Value = '  -1.83%  '
$ws.Range('D39').Value = '0.01840'
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('D40').Value = '2.726'
$ws.Range('E40').Value = '  -0.79%  '
$ws.Range('D41').Value = '0.9442'
$ws.Range('E41').Value = '  -1.46%  '
$ws.Range('D42').Value = '6.038'
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('D43').Value = '105.09'
$ws.Range('E43').Value = '  -2.45%  '
$ws.Range('D44').Value = '0.9985'
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('D45').Value = '1.983.84'
$ws.Range('E45').Value = '  -1.55%  '
$ws.Range('D46').Value = '0.5185'
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '64.80'
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.00000000122'
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.780'
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '9.466'
$ws.Range('E50').Value = '  -3.83%  '
$ws.Range('D51').Value = '0.07556'
$ws.Range('E51').Value = '  +13.66%  '
